# issue #5: add legislator_id, name, date into dataframe
# The stock-holdings sheet ("股票", 5th worksheet) gains three new trailing
# columns: date, legislator_name, legislator_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1): new column headers in H1:J1 -----------------------
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Match the bold/bordered header style used by the existing header cells
# (B1:G1) by copying their formatting onto the new header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Data row (row 2): new values in H2:J2 ----------------------------------
# H2 must stay a plain text string ("2013-11-12"), not get silently
# reinterpreted as a date serial, so force the cell to Text format before
# assigning it.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2013-11-12"
$ws.Range("I2").Value = "李慶華"
$ws.Range("J2").Value = 607

# Re-apply the plain data-row formatting (matching B2:G2) to the new data
# cells so they line up with the rest of row 2 instead of keeping the
# one-off "@" text format used to protect H2's value above.
$ws.Range("G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
